$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.247.70'
$ws.Range('E2').Value = '  +0.01%  '

$ws.Range('D3').Value = '3.564.50'
$ws.Range('E3').Value = '  +1.19%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '606.22'
$ws.Range('E5').Value = '  -0.07%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.35'
$ws.Range('E6').Value = '  -0.39%  '

$ws.Range('D7').Value = '3.563.04'
$ws.Range('E7').Value = '  +1.15%  '

$ws.Range('E9').Value = '  +2.58%  '

$ws.Range('E10').Value = '  -0.03%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.82'
$ws.Range('E11').Value = '  -2.90%  '

$ws.Range('E12').Value = '  -0.08%  '

$ws.Range('D13').Value = '4.169.79'
$ws.Range('E13').Value = '  +1.16%  '

$ws.Range('E14').Value = '  -0.22%  '

$ws.Range('E15').Value = '  -0.14%  '

$ws.Range('D16').Value = '3.571.91'
$ws.Range('E16').Value = '  +1.43%  '

$ws.Range('D17').Value = '66.291.81'
$ws.Range('E17').Value = '  -0.09%  '

$ws.Range('E18').Value = '  -0.59%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.42'
$ws.Range('E19').Value = '  +6.30%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.23'
$ws.Range('E20').Value = '  +0.26%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.80'
$ws.Range('E21').Value = '  -0.96%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '430.84'
$ws.Range('E22').Value = '  +1.12%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.613'
$ws.Range('E23').Value = '  +2.10%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '79.46'
$ws.Range('E24').Value = '  +1.47%  '

$ws.Range('D25').Value = '3.708.05'
$ws.Range('E25').Value = '  +1.35%  '

$ws.Range('E26').Value = '  -0.01%  '

$ws.Range('E27').Value = '  -2.00%  '

$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.50'
$ws.Range('E28').Value = '  +1.04%  '

$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.16'
$ws.Range('E29').Value = '  -1.37%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.90'
$ws.Range('E30').Value = '  -1.50%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.991'
$ws.Range('E31').Value = '  -0.83%  '

$ws.Range('D32').Value = '3.561.26'
$ws.Range('E32').Value = '  +1.47%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '25.43'
$ws.Range('E33').Value = '  +0.71%  '

$ws.Range('E34').Value = '  -1.71%  '

$ws.Range('E35').Value = '  -8.21%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.87'
$ws.Range('E36').Value = '  +1.25%  '

$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('E38').Value = '  -0.57%  '

$ws.Range('E39').Value = '  -0.07%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '173.84'
$ws.Range('E40').Value = '  +1.87%  '

$ws.Range('E41').Value = '  -0.83%  '

$ws.Range('E42').Value = '  +0.60%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.888'
$ws.Range('E43').Value = '  -0.21%  '

$ws.Range('E44').Value = '  +2.49%  '

$ws.Range('E45').Value = '  +1.08%  '

$ws.Range('E46').Value = '  -0.08%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.20'
$ws.Range('E47').Value = '  -1.22%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.45'
$ws.Range('E48').Value = '  +1.65%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '24.92'
$ws.Range('E49').Value = '  -3.60%  '

$ws.Range('E50').Value = '  -0.51%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '23.41'
$ws.Range('E51').Value = '  +4.50%  '
